$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("B2").Value = 0.00168175937904269
    $ws.Range("C2").Value = 0.999353169469599
    $ws.Range("E2").Value = 0.000258732212160414
    $ws.Range("F2").Value = 0.957309184993532
    $ws.Range("H2").Value = 0.979301423027167
    $ws.Range("I2").Value = 0.000129366106080207
    $ws.Range("K2").Value = 0.982147477360931
    $ws.Range("P2").Value = 0.999223803363519
    $ws.Range("Q2").Value = 0.0438551099611902
    $ws.Range("R2").Value = 0.000388098318240621
    $ws.Range("S2").Value = 0.000517464424320828
    $ws.Range("T2").Value = 0.986028460543338
    $ws.Range("U2").Value = 0.000646830530401035
    $ws.Range("W2").Value = 0.111384217335058
    $ws.Range("X2").Value = 0.000388098318240621
    $ws.Range("C3").Value = 0.000129366106080207
    $ws.Range("E3").Value = 0.000129366106080207
    $ws.Range("F3").Value = 0.000129366106080207
    $ws.Range("G3").Value = 0.983182406209573
    $ws.Range("H3").Value = 0.00168175937904269
    $ws.Range("K3").Value = 0.00646830530401035
    $ws.Range("L3").Value = 0.000517464424320828
    $ws.Range("N3").Value = 0.0939197930142303
    $ws.Range("Q3").Value = 0.000258732212160414
    $ws.Range("T3").Value = 0.00116429495472186
    $ws.Range("V3").Value = 0.000517464424320828
    $ws.Range("B4").Value = 0.998188874514877
    $ws.Range("C4").Value = 0.000517464424320828
    $ws.Range("E4").Value = 0.999611901681759
    $ws.Range("F4").Value = 0.0283311772315653
    $ws.Range("H4").Value = 0.0177231565329884
    $ws.Range("I4").Value = 0.99974126778784
    $ws.Range("K4").Value = 0.00827943078913325
    $ws.Range("M4").Value = 0.99987063389392
    $ws.Range("N4").Value = 0.000129366106080207
    $ws.Range("P4").Value = 0.000646830530401035
    $ws.Range("Q4").Value = 0.952263906856404
    $ws.Range("R4").Value = 0.999611901681759
    $ws.Range("S4").Value = 0.999353169469599
    $ws.Range("T4").Value = 0.0109961190168176
    $ws.Range("U4").Value = 0.999353169469599
    $ws.Range("W4").Value = 0.878783958602846
    $ws.Range("X4").Value = 0.999223803363519
    $ws.Range("B5").Value = 0.000129366106080207
    $ws.Range("F5").Value = 0.0134540750323415
    $ws.Range("G5").Value = 0.0166882276843467
    $ws.Range("H5").Value = 0.000129366106080207
    $ws.Range("I5").Value = 0.000129366106080207
    $ws.Range("K5").Value = 0.0018111254851229
    $ws.Range("L5").Value = 0.999353169469599
    $ws.Range("N5").Value = 0.901811125485123
    $ws.Range("T5").Value = 0.000776196636481242
    $ws.Range("V5").Value = 0.998965071151358
    $ws.Range("W5").Value = 0.000646830530401035
    $ws.Range("X5").Value = 0.000129366106080207
